$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 1041
$ws.Range("F4").Value = 5832
$ws.Range("F5").Value = 544
$ws.Range("F6").Value = 1034
$ws.Range("F7").Value = 1037
$ws.Range("F11").Value = 614
$ws.Range("F15").Value = 1991
$ws.Range("F16").Value = 1510
$ws.Range("F17").Value = 1073
$ws.Range("F20").Value = 399
$ws.Range("F21").Value = 634
$ws.Range("F22").Value = 227
$ws.Range("F23").Value = 1070
$ws.Range("F26").Value = 3500
$ws.Range("F27").Value = 193
$ws.Range("F28").Value = 133
$ws.Range("F30").Value = 158
$ws.Range("F32").Value = 492
$ws.Range("F37").Value = 322
$ws.Range("F38").Value = 825
$ws.Range("F39").Value = 104
$ws.Range("F40").Value = 70
$ws.Range("F41").Value = 79

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 36
$ws.Range("F4").Value = 665
$ws.Range("F6").Value = 352

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 1041
$ws.Range("F4").Value = 36
$ws.Range("F5").Value = 5832
$ws.Range("F6").Value = 544
$ws.Range("F7").Value = 1034
$ws.Range("F9").Value = 665
$ws.Range("F10").Value = 1037
$ws.Range("F13").Value = 353
$ws.Range("F16").Value = 614
$ws.Range("F21").Value = 1991
$ws.Range("F22").Value = 1510
$ws.Range("F23").Value = 1073
$ws.Range("F26").Value = 399
$ws.Range("F28").Value = 634
$ws.Range("F29").Value = 227
$ws.Range("F30").Value = 1070
$ws.Range("F31").Value = 3500
$ws.Range("F32").Value = 193
$ws.Range("F33").Value = 133
$ws.Range("F35").Value = 158
$ws.Range("F37").Value = 492
$ws.Range("F41").Value = 322
$ws.Range("F42").Value = 825
$ws.Range("F43").Value = 104
$ws.Range("F44").Value = 70
$ws.Range("F45").Value = 79
